$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Salario  - R$ 3200,00 `nNumero de dependentes - 2"
$ws.Range("B2").Value = "Salario liquido R$2.805,64`nINSS  R$352,00`nIR R$42,36`nFGTS R$ 256,00"

# Row 3
$enDash = [char]0x2013
$ws.Range("A3").Value = "Salario  - R$ 3200,00 `nNumero de dependentes $enDash 0"
$ws.Range("B3").Value = "Salario liquido R$2.775,60`nINSS R$352,00`nIR R$72,40`nFGTS R$ 256,00"

# Row 4
$ws.Range("A4").Value = "Salario  - R$ 1200,00 `nNumero de dependentes - 2"
$ws.Range("B4").Value = "Salario liquido R$2.775,60`nINSS R$352,00`nIR R$72,40`nFGTS R$ 256,00"

# Wrap text + row heights for the new rows
$ws.Range("A2:B4").WrapText = $true
$ws.Rows.Item(2).RowHeight = 46.45
$ws.Rows.Item(3).RowHeight = 46.45
$ws.Rows.Item(4).RowHeight = 46.45

# Final selection matches authored file
[void]$ws.Range("B4").Select()

# Tab-bar split ratio nudged slightly by the edit (cosmetic UI state)
$excel.ActiveWindow.TabRatio = 993
